$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the review row for nevilgreen12@gmail.com / vikicrestina@gmail.com
# (row 2), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# The row delete doesn't relocate the worksheet hyperlinks automatically,
# so drop the two stale ones and recreate the single hyperlink that
# survives on the new row 2 (the armonravid@gmail.com email, formerly on
# C3, now on C2).
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:armonravid@gmail.com", "", "", "armonravid@gmail.com")

# Hyperlinks.Add stamps the builtin "Hyperlink" style onto the cell;
# restore the original column formatting (copy it from the neighbouring
# D2 cell, which kept its formatting) instead of leaving the blue/underline
# look behind.
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to match the post-delete active cell.
$ws.Range("A2").Select()
